$d = $word.ActiveDocument

# Locate the "LOT2053 - ... (Requisito fraco)" run, including its trailing
# line break, inside the "Requisitos" bullet paragraph.
$rngSrc = $d.Content
$null = $rngSrc.Find.Execute("LOT2053*fraco)", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
$rngSrc.MoveEnd(1, 1)
$srcText = $rngSrc.Text

# Locate where the "LOT2007 - ... (Requisito fraco)" run starts - the
# LOT2053 entry needs to move in front of it.
$rngDst = $d.Content
$null = $rngDst.Find.Execute("LOT2007*fraco)", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
$dstStart = $rngDst.Start

# Insert a copy of the LOT2053 text (with its break) right before LOT2007.
$insertPoint = $d.Range($dstStart, $dstStart)
$insertPoint.InsertBefore($srcText)

# Remove the original LOT2053 occurrence, which now sits further along in
# the document because of the text we just inserted ahead of it.
$shift = $srcText.Length
$rngDel = $d.Range($rngSrc.Start + $shift, $rngSrc.End + $shift)
$rngDel.Delete()
